# Create Login Test Document - mark the task as Done and set its progress to 100%
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IT")

# Row 3 = "Create Login Test Document" task
# Column D = Status -> Done
$ws.Range("D3").Value = "Done"

# Column L = % (progress) -> 100
$ws.Range("L3").Value = 100

# Move the active selection, matching where the user left off
$ws.Range("F5").Select() | Out-Null
